# Trade #23 closed at 2026-02-17 20:53:47 - unknown UNKNOWN +0.000%
#
# - Summary sheet: refresh aggregate stats (capital, P&L, trade counts, win rate)
# - Strategy Status sheet: refresh MarketMaking strategy row
# - All Trades sheet: close out trade row 52 (early_exit) + append new open trade (row 85)
# - MarketMaking sheet: same two updates, mirrored into its own column layout (row 19 / row 52)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(3, 2).Value2  = 1400.4              # Current Capital
$summary.Cells.Item(4, 2).Value2  = 0.19                # Total P&L $
$summary.Cells.Item(5, 2).Value2  = 0.07000000000000001 # Total P&L %
$summary.Cells.Item(6, 2).Value2  = 51                  # Total Trades
$summary.Cells.Item(8, 2).Value2  = 21                  # Losing Trades
$summary.Cells.Item(9, 2).Value2  = 45.1                # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(5, 3).Value2 = 100.4   # Capital
$status.Cells.Item(5, 4).Value2 = 18      # Trades
$status.Cells.Item(5, 5).Value2 = 0.08    # P&L $
$status.Cells.Item(5, 6).Value2 = 0.4     # P&L %
$status.Cells.Item(5, 7).Value2 = 50      # Win Rate %

# ---------------------------------------------------------------------------
# All Trades
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out existing open trade on row 52 (early exit)
$allTrades.Cells.Item(52, 7).Value2  = 0.82          # Exit Price
$allTrades.Cells.Item(52, 8).Value2  = "CLOSED"      # Status
$allTrades.Cells.Item(52, 9).Value2  = -2.381        # P&L %
$allTrades.Cells.Item(52, 10).Value2 = -0.02         # P&L $
$allTrades.Cells.Item(52, 11).Value2 = 100.4         # Capital After
$allTrades.Cells.Item(52, 12).Value2 = "early_exit"  # Exit Reason
$allTrades.Cells.Item(52, 13).Value2 = 0.14          # Duration (min)

# Append new open trade as row 85 (keep the date/time as literal text, not
# an auto-converted date serial, by forcing a text number format first)
$allTrades.Cells.Item(85, 1).Value2  = 84
$allTrades.Cells.Item(85, 2).NumberFormat = "@"
$allTrades.Cells.Item(85, 2).Value2  = "2026-02-17"
$allTrades.Cells.Item(85, 2).Style = "Normal"
$allTrades.Cells.Item(85, 3).Value2  = "20:53:40"
$allTrades.Cells.Item(85, 4).Value2  = "MarketMaking"
$allTrades.Cells.Item(85, 5).Value2  = "DOWN"
$allTrades.Cells.Item(85, 6).Value2  = 0.84
$allTrades.Cells.Item(85, 8).Value2  = "OPEN"
$allTrades.Cells.Item(85, 9).Value2  = 0
$allTrades.Cells.Item(85, 10).Value2 = 0
$allTrades.Cells.Item(85, 11).Value2 = 100.4184370824165
$allTrades.Cells.Item(85, 13).Value2 = 0
$allTrades.Cells.Item(85, 14).Value2 = 0
$allTrades.Cells.Item(85, 15).Value2 = 0
$allTrades.Cells.Item(85, 16).Value2 = 0.6
$allTrades.Cells.Item(85, 17).Value2 = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking (per-strategy trade log; same two trades, different columns)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close out existing open trade on row 19 (early exit)
$mm.Cells.Item(19, 7).Value2  = 0.82          # Exit Price
$mm.Cells.Item(19, 8).Value2  = "CLOSED"      # Status
$mm.Cells.Item(19, 9).Value2  = -2.381        # P&L %
$mm.Cells.Item(19, 10).Value2 = -0.02         # P&L $
$mm.Cells.Item(19, 11).Value2 = 100.4         # Capital After
$mm.Cells.Item(19, 16).Value2 = "early_exit"  # Exit Reason
$mm.Cells.Item(19, 17).Value2 = 0.14          # Duration (min)

# Append new open trade as row 52
$mm.Cells.Item(52, 1).Value2  = 84
$mm.Cells.Item(52, 2).NumberFormat = "@"
$mm.Cells.Item(52, 2).Value2  = "2026-02-17"
$mm.Cells.Item(52, 2).Style = "Normal"
$mm.Cells.Item(52, 3).Value2  = "20:53:40"
$mm.Cells.Item(52, 4).Value2  = "MarketMaking"
$mm.Cells.Item(52, 5).Value2  = "DOWN"
$mm.Cells.Item(52, 6).Value2  = 0.84
$mm.Cells.Item(52, 8).Value2  = "OPEN"
$mm.Cells.Item(52, 9).Value2  = 0
$mm.Cells.Item(52, 10).Value2 = 0
$mm.Cells.Item(52, 11).Value2 = 100.4184370824165
$mm.Cells.Item(52, 12).Value2 = 0
$mm.Cells.Item(52, 13).Value2 = 0
$mm.Cells.Item(52, 14).Value2 = 0.6
$mm.Cells.Item(52, 15).Value2 = "Normal spread capture: 19600 bps"
$mm.Cells.Item(52, 17).Value2 = 0
